$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price (D) and 1h-volume-change (E) columns with the latest scrape.
# D-column values are plain text that often *look* numeric ("147.30", "1.00", ...).
# A bare Range.Value assignment lets Excel auto-convert those into real numbers
# (dropping significant trailing zeros / using a different internal type than the
# original inline-string cells), so every D-column write is routed through a
# Text-formatted scratch cell + Copy/PasteSpecial(values) round-trip, which keeps
# the literal text and leaves the destination cell style untouched.
function Set-TextValue($cellRef, $text) {
    $ws.Range("A52").NumberFormat = "@"
    $ws.Range("A52").Value = $text
    $ws.Range("A52").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

Set-TextValue "D2" "70.627.32"
$ws.Range("E2").Value = "  -2.55%  "
Set-TextValue "D3" "3.930.73"
$ws.Range("E3").Value = "  -2.74%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.25%  "
Set-TextValue "D5" "539.09"
$ws.Range("E5").Value = "  +4.13%  "
Set-TextValue "D6" "147.30"
$ws.Range("E6").Value = "  +0.03%  "
Set-TextValue "D7" "0.688"
$ws.Range("E7").Value = "  -3.42%  "
$ws.Range("E8").Value = "  +0.00%  "
Set-TextValue "D9" "0.733"
$ws.Range("E9").Value = "  -3.50%  "
$ws.Range("E10").Value = "  -4.30%  "
Set-TextValue "D11" "52.93"
$ws.Range("E11").Value = "  +13.57%  "
$ws.Range("E12").Value = "  -2.89%  "
Set-TextValue "D13" "10.48"
$ws.Range("E13").Value = "  -3.72%  "
Set-TextValue "D14" "4.543.34"
$ws.Range("E14").Value = "  -3.10%  "
Set-TextValue "D15" "3.938.36"
$ws.Range("E15").Value = "  -3.35%  "
Set-TextValue "D16" "13.86"
$ws.Range("E16").Value = "  -2.24%  "
Set-TextValue "D17" "20.22"
$ws.Range("E17").Value = "  -4.07%  "
$ws.Range("E18").Value = "  -0.78%  "
Set-TextValue "D19" "1.17"
$ws.Range("E19").Value = "  -3.83%  "
Set-TextValue "D20" "70.529.53"
$ws.Range("E20").Value = "  -2.56%  "
Set-TextValue "D21" "426.64"
$ws.Range("E21").Value = "  -3.58%  "
Set-TextValue "D22" "96.33"
$ws.Range("E22").Value = "  -7.04%  "
Set-TextValue "D23" "3.50"
$ws.Range("E23").Value = "  -1.50%  "
Set-TextValue "D24" "4.18"
$ws.Range("E24").Value = "  +5.18%  "
Set-TextValue "D25" "14.17"
$ws.Range("E25").Value = "  -3.27%  "
$ws.Range("E26").Value = "  -3.75%  "
Set-TextValue "D27" "10.48"
$ws.Range("E27").Value = "  -6.14%  "
$ws.Range("E28").Value = "  +0.85%  "
Set-TextValue "D29" "3.60"
$ws.Range("E29").Value = "  +16.21%  "
Set-TextValue "D30" "36.21"
$ws.Range("E30").Value = "  -4.59%  "
Set-TextValue "D31" "7.46"
$ws.Range("E31").Value = "  +9.36%  "
Set-TextValue "D32" "13.35"
$ws.Range("E32").Value = "  -2.36%  "
Set-TextValue "D33" "676.38"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("E34").Value = "  +0.39%  "
Set-TextValue "D35" "47.54"
$ws.Range("E35").Value = "  +13.62%  "
Set-TextValue "D36" "64.84"
$ws.Range("E36").Value = "  -4.11%  "
Set-TextValue "D37" "0.427"
$ws.Range("E37").Value = "  -0.99%  "
Set-TextValue "D38" "0.0₃0816"
$ws.Range("E38").Value = "  -5.12%  "
$ws.Range("E39").Value = "  -3.41%  "
Set-TextValue "D40" "0.148"
$ws.Range("E40").Value = "  -1.78%  "
Set-TextValue "D41" "1.00"
$ws.Range("E41").Value = "  +0.22%  "
Set-TextValue "D42" "3.30"
$ws.Range("E42").Value = "  +3.73%  "
Set-TextValue "D43" "0.998"
$ws.Range("E43").Value = "  -0.16%  "
Set-TextValue "D44" "0.0479"
$ws.Range("E44").Value = "  -2.36%  "
Set-TextValue "D45" "0.148"
$ws.Range("E45").Value = "  -4.65%  "
Set-TextValue "D46" "2.68"
$ws.Range("E46").Value = "  -1.98%  "
Set-TextValue "D47" "9.60"
$ws.Range("E47").Value = "  +4.80%  "
Set-TextValue "D48" "3.33"
$ws.Range("E48").Value = "  -4.50%  "
Set-TextValue "D49" "2.96"
$ws.Range("E49").Value = "  -3.63%  "
Set-TextValue "D50" "0.000271"
$ws.Range("E50").Value = "  +1.47%  "
Set-TextValue "D51" "145.08"
$ws.Range("E51").Value = "  +1.46%  "

$ws.Range("A52").ClearContents()
$ws.Range("A52").ClearFormats()
$ws.Rows("52:52").Delete()

